# MORLIST.xlsx - "Update to all model sheets"
#
# The "model" sheet's lookup table (A2:C9 - field name / type / flag) is
# re-sorted alphabetically by column A, and "model" becomes the active
# (selected) sheet/tab with cell E6 selected, replacing "settings" as the
# previously active tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")

# Re-sort A2:C9 ascending by column A. Using the worksheet Sort object (vs.
# Range.Sort) so the persisted sort state / sortState ref on the sheet is
# refreshed to cover the full current data range.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2"))
$ws.Sort.SetRange($ws.Range("A2:C9"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# "model" becomes the active sheet/tab (this clears tabSelected on whatever
# sheet was active before - "settings" - and sets it here), with E6 as the
# selected cell.
$ws.Activate()
$ws.Range("E6").Select()
